# Generate Report for Handback
# Row 7 ("4e223e4f-05ec-4d26-9767-c0f865b04092") has now received a
# handback that is not the latest version, so the per-language sheets
# (zh-cn / de-de) get their "Latest Target File" / "Latest Handback
# File" / "Latest Handback DateTime" / "Error Detail" columns populated,
# plus a new hyperlink on the "Latest Target File" cell (column I).

$wb = $excel.ActiveWorkbook

$mdName     = "4e223e4f-05ec-4d26-9767-c0f865b04092.md"
$targetUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/63105927f555732d512873ab38bef163ad865516/e2e/4e223e4f-05ec-4d26-9767-c0f865b04092.md"
$errorMsg   = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d2e8ecbfcb13c1997592e218fdadda9c947dbb5a/e2e/4e223e4f-05ec-4d26-9767-c0f865b04092.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/63105927f555732d512873ab38bef163ad865516/e2e/4e223e4f-05ec-4d26-9767-c0f865b04092.md."

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I7").Value = $mdName
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $targetUrl, "", "", $mdName)
$wsZh.Range("I7").Font.Underline = $true
$wsZh.Range("I7").Font.Color = $wsZh.Range("A7").Font.Color

$wsZh.Range("J7").Value = "4e223e4f-05ec-4d26-9767-c0f865b04092.76c04d3d4d00049a963fe8c5b5bf4bcdfdfa78c8.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-09-02 13:04:21"
$wsZh.Range("P7").Value = $errorMsg

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I7").Value = $mdName
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $targetUrl, "", "", $mdName)
$wsDe.Range("I7").Font.Underline = $true
$wsDe.Range("I7").Font.Color = $wsDe.Range("A7").Font.Color

$wsDe.Range("J7").Value = "4e223e4f-05ec-4d26-9767-c0f865b04092.76c04d3d4d00049a963fe8c5b5bf4bcdfdfa78c8.de-de.xlf"
$wsDe.Range("K7").Value = "2016-09-02 13:04:28"
$wsDe.Range("P7").Value = $errorMsg
